# Update column F ("dSF") values for several rows, per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -9
$ws.Range("F3").Value = -5
$ws.Range("F4").Value = 5
$ws.Range("F5").Value = -7
$ws.Range("F6").Value = 3
$ws.Range("F14").Value = 1
